$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

Set-TextCell $ws "D2" "65.903.25"
$ws.Range("E2").Value = "  +6.45%  "
Set-TextCell $ws "D3" "3.003.78"
$ws.Range("E3").Value = "  +3.74%  "
$ws.Range("E4").Value = "  +0.06%  "
Set-TextCell $ws "D5" "582.78"
$ws.Range("E5").Value = "  +2.63%  "
Set-TextCell $ws "D6" "161.41"
$ws.Range("E6").Value = "  +12.49%  "
$ws.Range("E7").Value = "  -0.03%  "
Set-TextCell $ws "D8" "2.999.75"
$ws.Range("E8").Value = "  +3.60%  "
$ws.Range("E9").Value = "  +2.90%  "
Set-TextCell $ws "D10" "7.01"
$ws.Range("E10").Value = "  +1.16%  "
Set-TextCell $ws "D11" "0.156"
$ws.Range("E11").Value = "  +6.32%  "
$ws.Range("E12").Value = "  +5.73%  "
$ws.Range("E13").Value = "  +8.26%  "
$ws.Range("E14").Value = "  +8.16%  "
$ws.Range("E15").Value = "  +0.75%  "
Set-TextCell $ws "D16" "65.898.87"
$ws.Range("E16").Value = "  +6.62%  "
Set-TextCell $ws "D17" "3.504.44"
$ws.Range("E17").Value = "  +3.83%  "
$ws.Range("E18").Value = "  +6.27%  "
Set-TextCell $ws "D19" "3.004.10"
$ws.Range("E19").Value = "  +3.83%  "
Set-TextCell $ws "D20" "457.35"
$ws.Range("E20").Value = "  +6.40%  "
Set-TextCell $ws "D21" "13.85"
$ws.Range("E21").Value = "  +6.90%  "
$ws.Range("E22").Value = "  +4.63%  "
Set-TextCell $ws "D23" "7.33"
$ws.Range("E23").Value = "  +6.70%  "
Set-TextCell $ws "D24" "82.21"
$ws.Range("E24").Value = "  +4.18%  "
Set-TextCell $ws "D25" "2.29"
$ws.Range("E25").Value = "  +13.10%  "
Set-TextCell $ws "D26" "12.34"
$ws.Range("E26").Value = "  +2.49%  "
Set-TextCell $ws "D27" "10.58"
$ws.Range("E27").Value = "  +5.12%  "
$ws.Range("E28").Value = "  -0.10%  "
Set-TextCell $ws "D29" "8.12"
$ws.Range("E29").Value = "  +15.55%  "
$ws.Range("E30").Value = "  +15.29%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell $ws "D31" "2.60"
$ws.Range("E31").Value = "  +4.06%  "
$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextCell $ws "D32" "0.0000103"
$ws.Range("E32").Value = "  -6.81%  "
Set-TextCell $ws "D33" "26.93"
$ws.Range("E33").Value = "  +5.14%  "
$ws.Range("E34").Value = "  +2.57%  "
$ws.Range("E35").Value = "  +0.02%  "
Set-TextCell $ws "D36" "0.991"
$ws.Range("E36").Value = "  +3.87%  "
Set-TextCell $ws "D37" "5.78"
$ws.Range("E37").Value = "  +7.18%  "
Set-TextCell $ws "D38" "2.14"
$ws.Range("E38").Value = "  +11.66%  "
Set-TextCell $ws "D39" "49.79"
$ws.Range("E39").Value = "  +1.81%  "
$ws.Range("E40").Value = "  +4.81%  "
Set-TextCell $ws "D41" "0.305"
$ws.Range("E41").Value = "  +13.60%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell $ws "D42" "0.122"
$ws.Range("E42").Value = "  +6.19%  "
$ws.Range("B43").Value = "Arweave"
$ws.Range("C43").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextCell $ws "D43" "43.75"
$ws.Range("E43").Value = "  +9.01%  "
Set-TextCell $ws "D44" "8.46"
$ws.Range("E44").Value = "  +4.10%  "
Set-TextCell $ws "D45" "382.39"
$ws.Range("E45").Value = "  +10.45%  "
$ws.Range("E46").Value = "  +5.69%  "
Set-TextCell $ws "D47" "2.788.37"
$ws.Range("E47").Value = "  +3.47%  "
Set-TextCell $ws "D48" "135.04"
$ws.Range("E48").Value = "  +2.61%  "
$ws.Range("E49").Value = "  -0.10%  "
Set-TextCell $ws "D50" "23.90"
$ws.Range("E50").Value = "  +10.41%  "
Set-TextCell $ws "D51" "0.107"
$ws.Range("E51").Value = "  +3.77%  "
